$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.192.03'
$ws.Range('E2').Value = '  +9.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.226.00'
$ws.Range('E3').Value = '  +3.71%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '397.68'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.37'
$ws.Range('E6').Value = '  +7.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.553'
$ws.Range('E7').Value = '  +2.51%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.619'
$ws.Range('E9').Value = '  +4.74%  '
$ws.Range('E10').Value = '  +5.62%  '
$ws.Range('E11').Value = '  +5.79%  '
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.742.53'
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.08'
$ws.Range('E14').Value = '  +2.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.03'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.215.81'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.05'
$ws.Range('E17').Value = '  +4.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.72'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '56.060.54'
$ws.Range('E19').Value = '  +8.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.31'
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('E21').Value = '  +5.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.00'
$ws.Range('E22').Value = '  +3.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '302.59'
$ws.Range('E23').Value = '  +13.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.43'
$ws.Range('E24').Value = '  +7.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.23'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.20'
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.17'
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.49'
$ws.Range('E28').Value = '  +3.48%  '
$ws.Range('E29').Value = '  +4.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.110'
$ws.Range('E31').Value = '  +2.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.14'
$ws.Range('E32').Value = '  +6.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0494'
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '36.37'
$ws.Range('E34').Value = '  +1.39%  '
$ws.Range('E35').Value = '  +2.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.35'
$ws.Range('E36').Value = '  +2.56%  '
$ws.Range('E37').Value = '  +25.21%  '
$ws.Range('E38').Value = '  +3.50%  '
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.92'
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '132.78'
$ws.Range('E41').Value = '  +2.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.02'
$ws.Range('E42').Value = '  +5.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.31'
$ws.Range('E43').Value = '  +4.04%  '
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.282'
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.20'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('E47').Value = '  +45.72%  '
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.132.88'
$ws.Range('E50').Value = '  +2.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0366'
$ws.Range('E51').Value = '  +9.80%  '

Write-Output "Updated cryptos list"
